$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.161.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.42%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.658.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("E4").Value = "  -0.47%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5226"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.66%  "

$ws.Range("E7").Value = "  -0.48%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2628"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06295"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07818"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.494"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.657.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.886.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5546"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8015"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.180.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.639"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "195.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.963"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.007"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1201"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.148"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.486"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05733"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.274"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.495"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.382"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.586"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9563"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.805"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.420"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5710"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01596"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.960"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.065.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8500"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.51%  "

$ws.Range("E43").Value = "  -0.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.796.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.29%  "

$ws.Range("E46").Value = "  +1.02%  "

$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4404"
$ws.Range("D48").Style = "Normal"

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.35%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₈102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05204"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.84%  "
